$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New case figures for the countries whose numbers changed in this refresh.
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes
$updates = @{
    "Banglades"              = @(368690,1125,281656,81686,0,23,5348)
    "Indonesia"              = @(303498,3992,228453,63894,0,96,11151)
    "Rumania"                = @(135900,1835,108135,22762,0,56,5003)
    "Oman"                   = @(101270,2685,90296,9997,0,42,977)
    "Emiratos Arabes Unidos" = @(98801,1041,88123,10252,0,0,426)
    "Suiza"                  = @(54384,0,45800,6507,0,1,2077)
    "Afganistan"             = @(39341,44,32852,5027,0,0,1462)
    "Malasia"                = @(12381,293,10283,1961,0,0,137)
    "Eslovenia"              = @(6498,168,4265,2078,0,0,155)
    "Hong Kong"              = @(5114,5,4861,148,0,0,105)
    "Bahamas"                = @(4409,77,2375,1938,0,0,96)
    "Sri Lanka"              = @(3396,1,3258,125,0,0,13)
    "Gibraltar"              = @(432,4,360,72,0,0,0)
}

$firstDataRow = 4
$lastDataRow = 220

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $country = $ws.Cells.Item($r, 1).Value()
    if ($updates.ContainsKey($country)) {
        $vals = $updates[$country]
        $ws.Cells.Item($r, 2).Value = $vals[0]
        $ws.Cells.Item($r, 3).Value = $vals[1]
        $ws.Cells.Item($r, 4).Value = $vals[2]
        $ws.Cells.Item($r, 5).Value = $vals[3]
        $ws.Cells.Item($r, 6).Value = $vals[4]
        $ws.Cells.Item($r, 7).Value = $vals[5]
        $ws.Cells.Item($r, 8).Value = $vals[6]
    }
}

# Re-rank the table: it is kept sorted by "Casos totales" (column B) descending.
$sortRange = $ws.Range("A$firstDataRow`:H$lastDataRow")
$sortKey = $ws.Range("B$firstDataRow`:B$lastDataRow")
$sortRange.Sort($sortKey, 2)

# Bump the "last updated" timestamp shown at the top of the sheet.
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 12:30"
